# This script updates the NATMI ligand-receptor pair statistics for the
# Jag2-Notch2 sheet with a refreshed TPM-based recomputation. For every
# data row (2-26) the ligand/receptor expression & specificity columns
# (G:J and M:P, which are shared across rows with the same sending /
# target cluster) and the edge-level weight & specificity columns
# (Q:T) are refreshed to the newly computed values, carrying forward
# any column values that did not change (K, L stay as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> full new G:T row values (columns G,H,I,J,K,L,M,N,O,P,Q,R,S,T)
$rowUpdates = @(
    @(2, @(12.01751833333333, 36.052555, 0.5714697501126412, 0.6009086862203552, 3, 1, 1.400501333333333, 4.201504, 0.00926314904242919, 0.009687730200823723, 16.83055044919111, 151.47495404272, 0.005293609468533161, 0.005821441227434241)),
    @(3, @(12.01751833333333, 36.052555, 0.5714697501126412, 0.6009086862203552, 3, 1, 21.254561, 63.763683, 0.1405812059498714, 0.1470248171880475, 255.4270764844516, 2298.843688360065, 0.08033790663470676, 0.0883484897382575)),
    @(4, @(12.01751833333333, 36.052555, 0.5714697501126412, 0.6009086862203552, 3, 1, 63.87756733333333, 191.632702, 0.4224968677952986, 0.4418622271050682, 767.6498365170678, 6908.84852865361, 0.2414441794623529, 0.2655188503801068)),
    @(5, @(12.01751833333333, 36.052555, 0.5714697501126412, 0.6009086862203552, 2, 1, 19.878555, 39.75711, 0.1314800731212866, 0.0916710195312133, 238.890899152675, 1433.34539491605, 0.07513688453141343, 0.05508591191098191)),
    @(6, @(12.01751833333333, 36.052555, 0.5714697501126412, 0.6009086862203552, 3, 1, 44.77944466666667, 134.338334, 0.2961787040911142, 0.3097542059748472, 538.1377972381523, 4843.24017514337, 0.169257170015635, 0.1861339929635747)),
    @(7, @(0.6575503333333333, 1.972651, 0.03126852934637925, 0.03287930968502149, 3, 1, 1.400501333333333, 4.201504, 0.00926314904242919, 0.009687730200823723, 0.9209001185671112, 8.288101067104, 0.0002896450476730819, 0.0003185258814178186)),
    @(8, @(0.6575503333333333, 1.972651, 0.03126852934637925, 0.03287930968502149, 3, 1, 21.254561, 63.763683, 0.1405812059498714, 0.1470248171880475, 13.97594367040367, 125.783493033633, 0.004395767563792939, 0.004834074495709484)),
    @(9, @(0.6575503333333333, 1.972651, 0.03126852934637925, 0.03287930968502149, 3, 1, 63.87756733333333, 191.632702, 0.4224968677952986, 0.4418622271050682, 42.00271569255577, 378.024441233002, 0.01321085570941061, 0.01452812500310084)),
    @(10, @(0.6575503333333333, 1.972651, 0.03126852934637925, 0.03287930968502149, 2, 1, 19.878555, 39.75711, 0.1314800731212866, 0.0916710195312133, 13.071150466435, 78.42690279861, 0.004111188524857038, 0.003014079840308416)),
    @(11, @(0.6575503333333333, 1.972651, 0.03126852934637925, 0.03287930968502149, 3, 1, 44.77944466666667, 134.338334, 0.2961787040911142, 0.3097542059748472, 29.44473876704823, 265.002648903434, 0.009261072500645581, 0.01018450446448493)),
    @(12, @(3.025265666666666, 9.075797, 0.1438606346668928, 0.1512715326742485, 3, 1, 1.400501333333333, 4.201504, 0.00926314904242919, 0.009687730200823723, 4.236888599854222, 38.131997398688, 0.001332602500257884, 0.00146547779561321)),
    @(13, @(3.025265666666666, 9.075797, 0.1438606346668928, 0.1512715326742485, 3, 1, 21.254561, 63.763683, 0.1405812059498714, 0.1470248171880475, 64.30069365337232, 578.7062428803509, 0.02022410151018566, 0.02224066943718713)),
    @(14, @(3.025265666666666, 9.075797, 0.1438606346668928, 0.1512715326742485, 3, 1, 63.87756733333333, 191.632702, 0.4224968677952986, 0.4418622271050682, 193.2466113237215, 1739.219501913494, 0.06078066754580595, 0.06684117632504054)),
    @(15, @(3.025265666666666, 9.075797, 0.1438606346668928, 0.1512715326742485, 2, 1, 19.878555, 39.75711, 0.1314800731212866, 0.0916710195312133, 60.13790994444499, 360.8274596666699, 0.01891480676527776, 0.01386721562629761)),
    @(16, @(3.025265666666666, 9.075797, 0.1438606346668928, 0.1512715326742485, 3, 1, 44.77944466666667, 134.338334, 0.2961787040911142, 0.3097542059748472, 135.4697165224665, 1219.227448702198, 0.04260845634536553, 0.04685699349010999)),
    @(17, @(3.090697, 6.181394, 0.146972094676554, 0.1030288518400537, 3, 1, 1.400501333333333, 4.201504, 0.00926314904242919, 0.009687730200823723, 4.328525269429334, 25.971151616576, 0.001361424418066933, 0.0009981157195270809)),
    @(18, @(3.090697, 6.181394, 0.146972094676554, 0.1030288518400537, 3, 1, 21.254561, 63.763683, 0.1405812059498714, 0.1470248171880475, 65.691407919017, 394.148447514102, 0.02066151431060863, 0.01514779810687832)),
    @(19, @(3.090697, 6.181394, 0.146972094676554, 0.1030288518400537, 3, 1, 63.87756733333333, 191.632702, 0.4224968677952986, 0.4418622271050682, 197.4262057244313, 1184.557234346588, 0.06209524965415815, 0.04552455793012423)),
    @(20, @(3.090697, 6.181394, 0.146972094676554, 0.1030288518400537, 2, 1, 19.878555, 39.75711, 0.1314800731212866, 0.0916710195312133, 61.438590302835, 245.75436121134, 0.01932390175486197, 0.009444759889308043)),
    @(21, @(3.090697, 6.181394, 0.146972094676554, 0.1030288518400537, 3, 1, 44.77944466666667, 134.338334, 0.2961787040911142, 0.3097542059748472, 138.3996952929327, 830.3981717575961, 0.04353000453885832, 0.031913620194216)),
    @(22, @(2.238110333333333, 6.714331, 0.1064289911975326, 0.1119116195803211, 3, 1, 1.400501333333333, 4.201504, 0.00926314904242919, 0.009687730200823723, 3.134476505980444, 28.210288553824, 0.0009858676078981291, 0.001084169576831372)),
    @(23, @(2.238110333333333, 6.714331, 0.1064289911975326, 0.1119116195803211, 3, 1, 21.254561, 63.763683, 0.1405812059498714, 0.1470248171880475, 47.57005260456366, 428.130473441073, 0.01496191593057738, 0.01645378541001502)),
    @(24, @(2.238110333333333, 6.714331, 0.1064289911975326, 0.1119116195803211, 3, 1, 63.87756733333333, 191.632702, 0.4224968677952986, 0.4418622271050682, 142.9650435169291, 1286.685391652362, 0.04496591542357094, 0.04944951746669585)),
    @(25, @(2.238110333333333, 6.714331, 0.1064289911975326, 0.1119116195803211, 2, 1, 19.878555, 39.75711, 0.1314800731212866, 0.0916710195312133, 44.49039935723499, 266.94239614341, 0.01399329154487635, 0.01025905226431733)),
    @(26, @(2.238110333333333, 6.714331, 0.1064289911975326, 0.1119116195803211, 3, 1, 44.77944466666667, 134.338334, 0.2961787040911142, 0.3097542059748472, 100.2213378293949, 901.992040464554, 0.03152200069060981, 0.03466509486246153))
)

foreach ($entry in $rowUpdates) {
    $rowNum = $entry[0]
    $vals = $entry[1]

    $arr = New-Object 'object[,]' 1,14
    for ($i = 0; $i -lt 14; $i++) {
        $arr[0, $i] = $vals[$i]
    }

    $rng = $ws.Range("G$rowNum" + ":T$rowNum")
    $rng.Value = $arr
}
